$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 2762.4375
$ws.Range("I51").Value = 2500
$ws.Range("J51").Value = 3339.8
$ws.Range("K51").Value = 2500
$ws.Range("L51").Value = 3339.8
$ws.Range("M51").Value = -2016
$ws.Range("N51").Value = -4307.8

$ws.Range("H111").Value = 13810.177
$ws.Range("I111").Value = 17252.54
$ws.Range("J111").Value = 2622.5
$ws.Range("K111").Value = 51757.62
$ws.Range("L111").Value = 7867.5
$ws.Range("M111").Value = -48690.62
$ws.Range("N111").Value = -14001.5

$ws.Range("H132").Value = 13415220
$ws.Range("I132").Value = 17598184
$ws.Range("K132").Value = 52794552
$ws.Range("M132").Value = -52792022

$ws.Range("H135").Value = 15157024
$ws.Range("I135").Value = 17549500
$ws.Range("K135").Value = 157945500
$ws.Range("M135").Value = -157942965

$ws.Range("H137").Value = 1880980.9
$ws.Range("I137").Value = 102300.4
$ws.Range("K137").Value = 306901.2
$ws.Range("M137").Value = -304351.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 258
$ws.Range("I25").Value = 258
$ws.Range("J25").Value = 0
$ws.Range("K25").Value = 258
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 144
$ws.Range("N25").ClearContents()

$ws.Range("H45").Value = 25014.818
$ws.Range("I45").Value = 17514.9
$ws.Range("K45").Value = 17514.9
$ws.Range("M45").Value = -17137.9

$ws.Range("H61").Value = 4626.857
$ws.Range("I61").Value = 4461.75
$ws.Range("K61").Value = 4461.75
$ws.Range("M61").Value = -4249.75

$ws.Range("H74").Value = 41681868
$ws.Range("I74").Value = 9301.25
$ws.Range("K74").Value = 9301.25
$ws.Range("M74").Value = -8427.25

$ws.Range("H77").Value = 41681868
$ws.Range("I77").Value = 9301.25
$ws.Range("K77").Value = 46506.25
$ws.Range("M77").Value = -42138.25

$ws.Range("H132").Value = 2790.8667
$ws.Range("I132").Value = 2125.25
$ws.Range("K132").Value = 6375.75
$ws.Range("M132").Value = -3845.75

$ws.Range("H135").Value = 45185
$ws.Range("J135").Value = 45185
$ws.Range("L135").Value = 45185
$ws.Range("N135").Value = -55325

$ws.Range("H136").Value = 4626.857
$ws.Range("I136").Value = 4461.75
$ws.Range("K136").Value = 13385.25
$ws.Range("M136").Value = -10835.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()

$ws.Range("H86").Value = 3747.842
$ws.Range("I86").Value = 2520.2
$ws.Range("K86").Value = 2520.2
$ws.Range("M86").Value = -1397.2

$ws.Range("H89").Value = 3747.842
$ws.Range("I89").Value = 2520.2
$ws.Range("K89").Value = 12601
$ws.Range("M89").Value = -6985

$ws.Range("H134").Value = 1498.4
$ws.Range("I134").Value = 1217.0312
$ws.Range("K134").Value = 3651.0936
$ws.Range("M134").Value = -1116.0936

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9187996
$ws.Range("I31").Value = 3958137.5
$ws.Range("J31").Value = 27783048
$ws.Range("K31").Value = 3958137.5
$ws.Range("L31").Value = 27783048
$ws.Range("M31").Value = -3957842.5
$ws.Range("N31").Value = -27783638

$ws.Range("H34").Value = 9187996
$ws.Range("I34").Value = 3958137.5
$ws.Range("J34").Value = 27783048
$ws.Range("K34").Value = 3958137.5
$ws.Range("L34").Value = 27783048
$ws.Range("M34").Value = -3957935.5
$ws.Range("N34").Value = -27783452

$ws.Range("H41").Value = 0
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("M41").ClearContents()

$ws.Range("H99").Value = 2956.4285
$ws.Range("I99").Value = 2671.818
$ws.Range("K99").Value = 2671.818
$ws.Range("M99").Value = -1173.818

$ws.Range("H126").Value = 2956.4285
$ws.Range("I126").Value = 2671.818
$ws.Range("K126").Value = 8015.454000000001
$ws.Range("M126").Value = -5545.454000000001

$ws.Range("H132").Value = 1859.7037
$ws.Range("I132").Value = 1742.8077
$ws.Range("K132").Value = 5228.4231
$ws.Range("M132").Value = -2698.4231

$ws.Range("H140").Value = 74650
$ws.Range("J140").Value = 74650
$ws.Range("L140").Value = 74650
$ws.Range("N140").Value = -85010

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 66740600
$ws.Range("I4").Value = 83425250
$ws.Range("K4").Value = 250275750
$ws.Range("M4").Value = -250275638

$ws.Range("H139").Value = 1973.8182
$ws.Range("I139").Value = 1746.2
$ws.Range("K139").Value = 5238.6
$ws.Range("M139").Value = -98.60000000000036

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 37076704
$ws.Range("I18").Value = 55565056
$ws.Range("K18").Value = 55565056
$ws.Range("M18").Value = -55564763

$ws.Range("H113").Value = 3588.9524
$ws.Range("I113").Value = 3260.6924
$ws.Range("J113").Value = 4122.375
$ws.Range("K113").Value = 3260.6924
$ws.Range("L113").Value = 4122.375
$ws.Range("M113").Value = -1090.6924
$ws.Range("N113").Value = -8462.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9910.223
$ws.Range("I7").Value = 10625.267
$ws.Range("K7").Value = 10625.267
$ws.Range("M7").Value = -10513.267

$ws.Range("H40").Value = 4635.05
$ws.Range("I40").Value = 3666.5
$ws.Range("K40").Value = 3666.5
$ws.Range("M40").Value = -3530.5

$ws.Range("H60").Value = 200000
$ws.Range("J60").Value = 200000
$ws.Range("L60").Value = 200000
$ws.Range("N60").Value = -201018

$ws.Range("H68").Value = 9999.5
$ws.Range("I68").Value = 9999.5
$ws.Range("J68").Value = 0
$ws.Range("K68").Value = 9999.5
$ws.Range("L68").Value = 0
$ws.Range("M68").Value = -9250.5
$ws.Range("N68").ClearContents()

$ws.Range("H71").Value = 9999.5
$ws.Range("I71").Value = 9999.5
$ws.Range("J71").Value = 0
$ws.Range("K71").Value = 49997.5
$ws.Range("L71").Value = 0
$ws.Range("M71").Value = -46253.5
$ws.Range("N71").ClearContents()

$ws.Range("H126").Value = 9910.223
$ws.Range("I126").Value = 10625.267
$ws.Range("K126").Value = 31875.801
$ws.Range("M126").Value = -29405.801

$ws.Range("H134").Value = 0
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("M134").ClearContents()
$ws.Range("N134").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1000
$ws.Range("J13").Value = 1000
$ws.Range("L13").Value = 1000
$ws.Range("N13").Value = -1280

$ws.Range("H49").Value = 22200
$ws.Range("J49").Value = 22200
$ws.Range("L49").Value = 22200
$ws.Range("N49").Value = -22660

$ws.Range("H138").Value = 90425
$ws.Range("J138").Value = 90425
$ws.Range("L138").Value = 90425
$ws.Range("N138").Value = -100705
Write-Host "Edit complete"
